# Rename worksheet tabs to survey codes
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "TEST_IMPORT_SURVEY_RESP_1_test"
$wb.Worksheets.Item(2).Name = "TEST_IMPORT_SURVEY_RESP_2_test"

# Scroll the visible sheet tabs so the second sheet is the first one shown
# (bookViews/workbookView firstSheet="1")
$wb.Windows.Item(1).ScrollWorkbookTabs(1, 2)
